# WALCL.xlsx update
# 1) "Data" sheet: insert 5 new (older) weekly observations at the top
#    (rows 2-6), shifting the existing series down, and append 2 new
#    (newer) weekly observations at the bottom of the series.
# 2) "SeriesInfo" sheet: refresh the FRED metadata fields that changed
#    (realtime_start, realtime_end, observation_end, last_updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Insert 5 new rows right after the header row, pushing the existing
#     data (currently rows 2..113) down to rows 7..118.
$ws.Range("A2:A6").EntireRow.Insert()

$newFront = @(
    @(44440, 8349.173000000001),
    @(44447, 8357.314),
    @(44454, 8448.77),
    @(44461, 8489.824000000001),
    @(44468, 8447.981)
)

$r = 2
foreach ($pair in $newFront) {
    $ws.Range("A$r").Value = $pair[0]
    $ws.Range("B$r").Value = $pair[1]
    $r = $r + 1
}

# Match the formatting used by the rest of the date / value columns
# (A column uses the bordered/centered date style, B column is plain).
$ws.Range("A7").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B2:B6").PasteSpecial(-4122)

# --- Append 2 new rows at the end of the series (old last row was 113,
#     now sits at row 118 after the insert above, so new rows are 119/120).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$newBack = @(
    @(45259, 7796.145),
    @(45266, 7737.385)
)

foreach ($pair in $newBack) {
    $lastRow = $lastRow + 1
    $ws.Range("A$lastRow").Value = $pair[0]
    $ws.Range("B$lastRow").Value = $pair[1]
}

$ws.Range("A118").Copy()
$ws.Range("A119:A120").PasteSpecial(-4122)
$ws.Range("B118").Copy()
$ws.Range("B119:B120").PasteSpecial(-4122)

# --- SeriesInfo sheet metadata refresh
$ws2 = $wb.Worksheets.Item("SeriesInfo")

# These look like dates, so force them to stay plain text (leading
# apostrophe), then strip the resulting quote-prefix style so the cell
# matches the unstyled text cells around it.
$ws2.Range("B3").Value = "'2023-12-08"
$ws2.Range("B4").Value = "'2023-12-08"
$ws2.Range("B7").Value = "'2023-12-06"

$ws2.Range("B5").Copy()
$ws2.Range("B3").PasteSpecial(-4122)
$ws2.Range("B4").PasteSpecial(-4122)
$ws2.Range("B7").PasteSpecial(-4122)

# This one already fails Excel's date auto-detection (has a UTC offset
# suffix) so it stays plain text without any extra handling.
$ws2.Range("B14").Value = "2023-12-07 15:34:03-06"
